# "Add files via upload" - the author re-uploaded a newer version of the
# "by Coach" roster sheet. Comparing the two versions, the only meaningful
# content change is that the "Started" (Yes/No) flag was toggled for eight
# players (rows 14, 15, 18, 21, 45, 46, 83, 84). Everything else in the raw
# diff (random save GUIDs, the author's local absolute file path, the Excel
# window's pixel geometry, and sub-pixel column/row metrics coming from a
# different machine's font rendering) is incidental noise produced by
# Excel/Windows on save and not a deliberate edit, so we just reproduce the
# actual data change here.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Trevon Brazile: No -> Yes
$ws.Range("C14").Value = "Yes"

# Ja'Kobi Gillespie: No -> Yes
$ws.Range("C15").Value = "Yes"

# Malik Dia: Yes -> No
$ws.Range("C18").Value = "No"

# Meleek Thomas: Yes -> No
$ws.Range("C21").Value = "No"

# Boogie Fland: No -> Yes
$ws.Range("C45").Value = "Yes"

# AJ Storr: Yes -> No
$ws.Range("C46").Value = "No"

# Max Mackinnon: No -> Yes
$ws.Range("C83").Value = "Yes"

# Tahaad Pettiford: Yes -> No
$ws.Range("C84").Value = "No"
